$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (index values for columns B:E) were re-keyed to 16/20/16/20
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) — corrected/replacement measurements for columns B:E
$ws.Range("B2").Value = 12.232762440473394
$ws.Range("C2").Value = 11.559030914856949
$ws.Range("D2").Value = 13.066501868432459
$ws.Range("E2").Value = 12.322849670143519

# Row 3 (STR) — corrected/replacement measurements for columns B:E
$ws.Range("B3").Value = 11.096057042646507
$ws.Range("C3").Value = 10.465699751611059
$ws.Range("D3").Value = 12.623249442862919
$ws.Range("E3").Value = 11.561693571943357

# Selection now only spans the edited columns instead of the whole used range
$ws.Range("B1:E3").Select()
